$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row 2
$ws1.Range("D2").Value = 1
$ws1.Range("H2").Value = 4.49
$ws1.Range("L2").Value = 1

# Row 3
$ws1.Range("D3").Value = 1
$ws1.Range("H3").Value = 2.96
$ws1.Range("L3").Value = 0.9

# Row 4
$ws1.Range("D4").Value = 1
$ws1.Range("H4").Value = 1.82
$ws1.Range("I4").Value = "Low"
$ws1.Range("J4").Value = "Normal"
$ws1.Range("L4").Value = 1.18

# Row 5
$ws1.Range("D5").Value = 2
$ws1.Range("H5").Value = 0.5600000000000001
$ws1.Range("I5").Value = "Low"
$ws1.Range("L5").Value = 1.17

# Row 6
$ws1.Range("L6").Value = 1.06

# Row 7
$ws1.Range("D7").Value = 1
$ws1.Range("L7").Value = 1

# Row 8
$ws1.Range("D8").Value = 1
$ws1.Range("L8").Value = 1.16

# Row 9
$ws1.Range("D9").Value = 1
$ws1.Range("L9").Value = 1.11

# Row 10
$ws1.Range("L10").Value = 0.9

# Row 11
$ws1.Range("D11").Value = 3
$ws1.Range("L11").Value = 0.88

# Row 12
$ws1.Range("D12").Value = 1
$ws1.Range("L12").Value = 0.92

# Row 13
$ws1.Range("D13").Value = 1
$ws1.Range("L13").Value = 0.83

# Row 14
$ws1.Range("D14").Value = 1
$ws1.Range("L14").Value = 0.8

# Row 15
$ws1.Range("D15").Value = 1
$ws1.Range("L15").Value = 1.09

# Row 16
$ws1.Range("L16").Value = 0.9399999999999999

# Row 17
$ws1.Range("D17").Value = 1
$ws1.Range("L17").Value = 0.84

# --- Sheet: Summary ---
# These cells hold numeric-looking text (e.g. "49"), so force the
# number format to Text first; otherwise Excel auto-converts the
# assigned string into a real number.
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("B9:B11").NumberFormat = "@"
$ws2.Range("B9").Value = "36"
$ws2.Range("B10").Value = "18"
$ws2.Range("B11").Value = "8"

Write-Output "edit applied"
